# Applies the post-processed PoS results update described in the commit:
# "All results, tables, etc updated after changing Indonesian PoS and PUD tests data"
# The workbook stores only literal values (no formulas), so the refreshed metrics
# are written directly. A handful of cells also carry a manually-applied
# bold+underline "row max" highlight that moves to a different column when the
# recomputed values change which entry is the row maximum; those are swapped via
# a Copy/PasteSpecial(Formats) round-trip through a scratch cell so the existing
# shared cell style is reused instead of a new one being minted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Updated metric values ----
$ws.Range("C2").Value = 0.9824459770114943
$ws.Range("D2").Value = 0.8229394636015326
$ws.Range("E2").Value = 0.7516291187739463
$ws.Range("F2").Value = 0.8470927203065134
$ws.Range("G2").Value = 0.7861885057471264
$ws.Range("H2").Value = 0.5591141762452108
$ws.Range("I2").Value = 0.574663601532567
$ws.Range("L2").Value = 0.7238620689655172
$ws.Range("M2").Value = 0.667656704980843
$ws.Range("N2").Value = 0.6131157088122605
$ws.Range("O2").Value = 0.5895019157088123
$ws.Range("P2").Value = 0.5680674329501916
$ws.Range("Q2").Value = 0.6250727969348659
$ws.Range("R2").Value = 0.6401900383141762
$ws.Range("S2").Value = 0.7131310344827586
$ws.Range("T2").Value = 0.6425164750957855
$ws.Range("U2").Value = 0.5379708812260536
$ws.Range("C3").Value = 0.844756728092271
$ws.Range("D3").Value = 0.9882388591125169
$ws.Range("E3").Value = 0.777674665252552
$ws.Range("F3").Value = 0.8765174901990492
$ws.Range("G3").Value = 0.7590954716766728
$ws.Range("H3").Value = 0.5635499327664248
$ws.Range("I3").Value = 0.5846858961004526
$ws.Range("L3").Value = 0.8137535273953145
$ws.Range("M3").Value = 0.6331319482585557
$ws.Range("N3").Value = 0.5701785951023655
$ws.Range("O3").Value = 0.5201984810893733
$ws.Range("P3").Value = 0.5320732561883298
$ws.Range("Q3").Value = 0.5775458798128823
$ws.Range("R3").Value = 0.6349311566068825
$ws.Range("S3").Value = 0.7232628169921024
$ws.Range("T3").Value = 0.6829794890248291
$ws.Range("U3").Value = 0.5205772617942842
$ws.Range("L4").Value = 0.7793214614676083
$ws.Range("L5").Value = 0.8012414069278516
$ws.Range("L6").Value = 0.6810456308886015
$ws.Range("C7").Value = 0.5141525141525142
$ws.Range("D7").Value = 0.5670163170163171
$ws.Range("E7").Value = 0.4023476523476524
$ws.Range("F7").Value = 0.6224608724608724
$ws.Range("G7").Value = 0.5213120213120214
$ws.Range("H7").Value = 0.9492174492174492
$ws.Range("I7").Value = 0.5404595404595405
$ws.Range("L7").Value = 0.5937395937395937
$ws.Range("M7").Value = 0.5826673326673326
$ws.Range("N7").Value = 0.5621878121878122
$ws.Range("O7").Value = 0.5382950382950383
$ws.Range("P7").Value = 0.5938228438228438
$ws.Range("Q7").Value = 0.5105727605727606
$ws.Range("R7").Value = 0.4392274392274392
$ws.Range("S7").Value = 0.5322177822177823
$ws.Range("T7").Value = 0.491008991008991
$ws.Range("U7").Value = 0.4234099234099234
$ws.Range("L8").Value = 0.6271016311166876
$ws.Range("L9").Value = 0.3885852522175433
$ws.Range("L10").Value = 0.4925276620204052
$ws.Range("C11").Value = 0.7381154499151104
$ws.Range("D11").Value = 0.8202886247877759
$ws.Range("E11").Value = 0.750169779286927
$ws.Range("F11").Value = 0.8225806451612904
$ws.Range("G11").Value = 0.6646010186757215
$ws.Range("H11").Value = 0.6289473684210526
$ws.Range("I11").Value = 0.580730050933786
$ws.Range("L11").Value = 0.9275891341256367
$ws.Range("M11").Value = 0.7298811544991511
$ws.Range("N11").Value = 0.6679966044142615
$ws.Range("O11").Value = 0.5992359932088285
$ws.Range("P11").Value = 0.5795415959252971
$ws.Range("Q11").Value = 0.6387945670628183
$ws.Range("R11").Value = 0.6155348047538201
$ws.Range("S11").Value = 0.7343803056027165
$ws.Range("T11").Value = 0.666723259762309
$ws.Range("U11").Value = 0.5791171477079796
$ws.Range("C12").Value = 0.7401306447032093
$ws.Range("D12").Value = 0.8178074410678784
$ws.Range("E12").Value = 0.739136608917921
$ws.Range("F12").Value = 0.7874183470604942
$ws.Range("G12").Value = 0.7884597178831771
$ws.Range("H12").Value = 0.6336741456025751
$ws.Range("I12").Value = 0.5858657578339487
$ws.Range("L12").Value = 0.7567452428287418
$ws.Range("M12").Value = 0.9651614124775159
$ws.Range("N12").Value = 0.7496923222569346
$ws.Range("O12").Value = 0.6361355675470983
$ws.Range("P12").Value = 0.6616964877402253
$ws.Range("Q12").Value = 0.7340244248792956
$ws.Range("R12").Value = 0.636324907696677
$ws.Range("S12").Value = 0.7894537536684654
$ws.Range("T12").Value = 0.6491527028306352
$ws.Range("U12").Value = 0.6231657673009562
$ws.Range("L13").Value = 0.6954951998030688
$ws.Range("L14").Value = 0.568638510893323
$ws.Range("C15").Value = 0.3458646616541353
$ws.Range("D15").Value = 0.4875709682369188
$ws.Range("E15").Value = 0.3107257940770293
$ws.Range("F15").Value = 0.3551480742673009
$ws.Range("G15").Value = 0.3672702163572196
$ws.Range("H15").Value = 0.4797452815712751
$ws.Range("I15").Value = 0.4478287555623753
$ws.Range("L15").Value = 0.4375479515114317
$ws.Range("M15").Value = 0.4686205309191346
$ws.Range("N15").Value = 0.5049102347705999
$ws.Range("O15").Value = 0.4796685591529845
$ws.Range("P15").Value = 0.9703851465398189
$ws.Range("Q15").Value = 0.518796992481203
$ws.Range("R15").Value = 0.3775510204081632
$ws.Range("S15").Value = 0.4825840110480282
$ws.Range("T15").Value = 0.4502838729476753
$ws.Range("U15").Value = 0.3609022556390977
$ws.Range("L16").Value = 0.6762239126194656
$ws.Range("L17").Value = 0.6526362885462555
$ws.Range("L18").Value = 0.5704564973244368
$ws.Range("L19").Value = 0.3325822442541685
$ws.Range("L20").Value = 0.3127427074866793
$ws.Range("C22").Value = 0.5722987292231136
$ws.Range("D22").Value = 0.6305840659775428
$ws.Range("E22").Value = 0.5652665511173515
$ws.Range("F22").Value = 0.585066026977656
$ws.Range("G22").Value = 0.5870086938910211
$ws.Range("H22").Value = 0.5228748656994235
$ws.Range("I22").Value = 0.5177809287672387
$ws.Range("L22").Value = 0.6057914883337052
$ws.Range("M22").Value = 0.577684774271471
$ws.Range("N22").Value = 0.5511509600906036
$ws.Range("O22").Value = 0.5062225230588137
$ws.Range("P22").Value = 0.5240065938287732
$ws.Range("Q22").Value = 0.5426777771897285
$ws.Range("R22").Value = 0.5232566590321321
$ws.Range("S22").Value = 0.6255520324174779
$ws.Range("T22").Value = 0.5787243939322209
$ws.Range("U22").Value = 0.5147575270289987
$ws.Range("E27").Value = 0.7606545228767788
$ws.Range("F27").Value = 0.8349955452746061
$ws.Range("G27").Value = 0.7910473794796405
$ws.Range("H27").Value = 0.5822570184340591
$ws.Range("I27").Value = 0.5662236839428239
$ws.Range("L27").Value = 0.7598448191289786
$ws.Range("M27").Value = 0.6964224835308853
$ws.Range("N27").Value = 0.6196752488483311
$ws.Range("O27").Value = 0.5607685556429474
$ws.Range("P27").Value = 0.5645635229037982
$ws.Range("Q27").Value = 0.6247274572077183
$ws.Range("R27").Value = 0.6632856326629957
$ws.Range("S27").Value = 0.7262954547954331
$ws.Range("T27").Value = 0.6606699513757367
$ws.Range("U27").Value = 0.5835648231312077
$ws.Range("C28").Value = 0.4985774581922485
$ws.Range("D28").Value = 0.5700334119386324
$ws.Range("E28").Value = 0.5059049721928737
$ws.Range("F28").Value = 0.5446387450836305
$ws.Range("G28").Value = 0.5226400604572229
$ws.Range("H28").Value = 0.5900117347391555
$ws.Range("I28").Value = 0.533030527018032
$ws.Range("L28").Value = 0.5254885347735574
$ws.Range("M28").Value = 0.5590017005129541
$ws.Range("N28").Value = 0.5422048118609036
$ws.Range("O28").Value = 0.5035731314559063
$ws.Range("P28").Value = 0.5215265151431913
$ws.Range("Q28").Value = 0.5200398971652362
$ws.Range("R28").Value = 0.4757914465278972
$ws.Range("S28").Value = 0.5882018015538362
$ws.Range("T28").Value = 0.5356578790766559
$ws.Range("U28").Value = 0.4617005043044403
$ws.Range("C29").Value = 0.5796113438562092
$ws.Range("E29").Value = 0.5652504735135506
$ws.Range("F29").Value = 0.5902810511217063
$ws.Range("G29").Value = 0.5937407781306955
$ws.Range("H29").Value = 0.556149996550704
$ws.Range("I29").Value = 0.5398011411078324
$ws.Range("L29").Value = 0.6269301635312061
$ws.Range("M29").Value = 0.6109757188274048
$ws.Range("N29").Value = 0.6310170306246662
$ws.Range("O29").Value = 0.5800759149608058
$ws.Range("P29").Value = 0.619271167529173
$ws.Range("Q29").Value = 0.6217717019798843
$ws.Range("R29").Value = 0.5207570173825464
$ws.Range("T29").Value = 0.5702462538921501
$ws.Range("U29").Value = 0.5261375775177617
$ws.Range("L30").Value = 0.467104434402885
$ws.Range("C32").Value = 0.5764492702479743
$ws.Range("D32").Value = 0.6338476432389293
$ws.Range("E32").Value = 0.5689776546503686
$ws.Range("F32").Value = 0.5872667930870294
$ws.Range("G32").Value = 0.5906109782156538
$ws.Range("H32").Value = 0.5170837854753038
$ws.Range("I32").Value = 0.5133769932974774
$ws.Range("L32").Value = 0.5948419879591568
$ws.Range("M32").Value = 0.571431359552665
$ws.Range("N32").Value = 0.5474273263076017
$ws.Range("O32").Value = 0.502978982997487
$ws.Range("P32").Value = 0.5216267906794332
$ws.Range("Q32").Value = 0.538964539690135
$ws.Range("R32").Value = 0.5082459724526249
$ws.Range("S32").Value = 0.6124797180365676
$ws.Range("T32").Value = 0.5736575768409075
$ws.Range("U32").Value = 0.5102358063845964
$ws.Range("C37").Value = 0.8118086179250156
$ws.Range("D37").Value = 0.6361085071686206
$ws.Range("E37").Value = 0.6132314536267361
$ws.Range("F37").Value = 0.6584539654913433
$ws.Range("C38").Value = 0.5283589295729215
$ws.Range("D38").Value = 0.5495102655102483
$ws.Range("E38").Value = 0.5292692112276383
$ws.Range("F38").Value = 0.5153379078657074
$ws.Range("C39").Value = 0.593560226091286
$ws.Range("D39").Value = 0.5742937670632475
$ws.Range("E39").Value = 0.6126223067843868
$ws.Range("F39").Value = 0.5644333939952713
$ws.Range("D40").Value = 0.4071578159004677
$ws.Range("C42").Value = 0.5179710845423162
$ws.Range("D42").Value = 0.5391866967107786
$ws.Range("E42").Value = 0.5111069641991569
$ws.Range("F42").Value = 0.579408422450774

# ---- Row-max highlight (bold+underline) swaps ----
# The workbook marks each row's maximum score with a bold+underlined cell style
# (shared style index 8 vs the plain style index 9). Recomputing the table moved
# the row maximum to a different column in two rows, so swap the formatting between
# the old and new "winner" cells via a scratch cell (Z1, cleared afterwards) rather
# than setting Font properties directly, so the existing shared styles are reused.
$scratch = $ws.Range("Z1")

# Row 27: row max moved from D27 to C27
$c27 = $ws.Range("C27")
$d27 = $ws.Range("D27")
$scratch.Value = 0
$c27.Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null
$d27.Copy() | Out-Null
$c27.PasteSpecial(-4122) | Out-Null
$scratch.Copy() | Out-Null
$d27.PasteSpecial(-4122) | Out-Null
$c27.Value = 0.8380891922578954
$d27.Value = 0.8342564497361569
$scratch.Clear() | Out-Null

# Row 29: row max moved from D29 to S29
$d29 = $ws.Range("D29")
$s29 = $ws.Range("S29")
$scratch.Value = 0
$d29.Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null
$s29.Copy() | Out-Null
$d29.PasteSpecial(-4122) | Out-Null
$scratch.Copy() | Out-Null
$s29.PasteSpecial(-4122) | Out-Null
$d29.Value = 0.6389174838342686
$s29.Value = 0.6405927271886271
$scratch.Clear() | Out-Null
